{"js": "// Update the date line and the 25 two-digit multiplication problems to the\n// new \"output generated at c986bee\" values.\nconst replacements = [\n  [\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"],\n  [\"58\u00d759=\", \"89\u00d766=\"],\n  [\"44\u00d765=\", \"71\u00d780=\"],\n  [\"50\u00d753=\", \"82\u00d733=\"],\n  [\"35\u00d765=\", \"23\u00d776=\"],\n  [\"41\u00d732=\", \"41\u00d745=\"],\n  [\"77\u00d734=\", \"70\u00d741=\"],\n  [\"17\u00d741=\", \"15\u00d711=\"],\n  [\"99\u00d776=\", \"98\u00d752=\"],\n  [\"90\u00d776=\", \"33\u00d765=\"],\n  [\"67\u00d794=\", \"73\u00d726=\"],\n  [\"56\u00d754=\", \"22\u00d712=\"],\n  [\"50\u00d786=\", \"61\u00d737=\"],\n  [\"11\u00d759=\", \"13\u00d766=\"],\n  [\"59\u00d750=\", \"76\u00d790=\"],\n  [\"30\u00d760=\", \"65\u00d768=\"],\n  [\"42\u00d753=\", \"97\u00d759=\"],\n  [\"70\u00d723=\", \"89\u00d773=\"],\n  [\"51\u00d750=\", \"70\u00d792=\"],\n  [\"41\u00d734=\", \"42\u00d772=\"],\n  [\"87\u00d789=\", \"16\u00d795=\"],\n  [\"50\u00d779=\", \"97\u00d750=\"],\n  [\"22\u00d785=\", \"87\u00d767=\"],\n  [\"28\u00d730=\", \"74\u00d752=\"],\n  [\"21\u00d715=\", \"50\u00d753=\"],\n  [\"95\u00d791=\", \"53\u00d788=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 two-digit multiplication problems to the\n# new \"output generated at c986bee\" values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"),\n    @(\"58\u00d759=\", \"89\u00d766=\"),\n    @(\"44\u00d765=\", \"71\u00d780=\"),\n    @(\"50\u00d753=\", \"82\u00d733=\"),\n    @(\"35\u00d765=\", \"23\u00d776=\"),\n    @(\"41\u00d732=\", \"41\u00d745=\"),\n    @(\"77\u00d734=\", \"70\u00d741=\"),\n    @(\"17\u00d741=\", \"15\u00d711=\"),\n    @(\"99\u00d776=\", \"98\u00d752=\"),\n    @(\"90\u00d776=\", \"33\u00d765=\"),\n    @(\"67\u00d794=\", \"73\u00d726=\"),\n    @(\"56\u00d754=\", \"22\u00d712=\"),\n    @(\"50\u00d786=\", \"61\u00d737=\"),\n    @(\"11\u00d759=\", \"13\u00d766=\"),\n    @(\"59\u00d750=\", \"76\u00d790=\"),\n    @(\"30\u00d760=\", \"65\u00d768=\"),\n    @(\"42\u00d753=\", \"97\u00d759=\"),\n    @(\"70\u00d723=\", \"89\u00d773=\"),\n    @(\"51\u00d750=\", \"70\u00d792=\"),\n    @(\"41\u00d734=\", \"42\u00d772=\"),\n    @(\"87\u00d789=\", \"16\u00d795=\"),\n    @(\"50\u00d779=\", \"97\u00d750=\"),\n    @(\"22\u00d785=\", \"87\u00d767=\"),\n    @(\"28\u00d730=\", \"74\u00d752=\"),\n    @(\"21\u00d715=\", \"50\u00d753=\"),\n    @(\"95\u00d791=\", \"53\u00d788=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n}\n"}
